# Faculty/major seeder data update (src/seeders/files/xlsx/faculties.xlsx).
# Appends the full faculty list under the existing "GENERAL" row, matching
# the committed worksheet content exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# code, name pairs for rows 3..15 (row 2 "GENERAL"/"General" already exists).
$facultyData = @(
    ,@("FIKTI", "Fakultas Ilmu Komputer dan Teknologi Informasi")
    ,@("FIKF", "Fakultas Ilmu Kesehatan dan Farmasi")
    ,@("FK", "Fakultas Kedokteran")
    ,@("FTSP", "Fakultas Teknik Sipil dan Perencanaan")
    ,@("FE", "Fakultas Ekonomi")
    ,@("FPSI", "Fakultas Psikologi")
    ,@("FTI", "Fakultas Teknologi Industri")
    ,@("FSB", "Fakultas Sastra dan Budaya")
    ,@("FIKOM", "Fakultas Ilmu Komunikasi")
    ,@("FTI_DIPLOMA", "Fakultas Teknologi Informasi")
    ,@("FBK", "Fakultas Bisnis dan Kewirausahaan")
    ,@("MAGISTER", "Program Magister")
    ,@("PROFESI", "Program Profesi")
)

# Row 3 (FIKTI) registers its "name" (col B) shared string before its
# "code" (col A) shared string, so set B before A there to reproduce the
# original authoring order; every other row sets A before B.
$row = 3
$first = $true
foreach ($item in $facultyData) {
    if ($first) {
        $ws.Cells.Item($row, 2).Value = $item[1]
        $ws.Cells.Item($row, 1).Value = $item[0]
        $first = $false
    } else {
        $ws.Cells.Item($row, 1).Value = $item[0]
        $ws.Cells.Item($row, 2).Value = $item[1]
    }
    $ws.Cells.Item($row, 3).Value = $true
    $row++
}

# Column A/B were widened to fit the new (longer) codes/names.
$ws.Columns.Item(1).ColumnWidth = 11.67
$ws.Columns.Item(2).ColumnWidth = 42.67

# Move the active selection, matching the saved worksheet view.
$ws.Range("E7").Select() | Out-Null
